# Update "想去人数" (want-to-go count) values in column F for rows 3, 4, 5, 7
# on both the "展览" and "全部类型" worksheets, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    3 = 357
    4 = 1154
    5 = 1185
    7 = 5986
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
